$d = $word.ActiveDocument

$d.Content.Find.Execute("441×2=882", $true, $false, $false, $false, $false, $true, 1, $false, "426×4=1704", 2) | Out-Null
$d.Content.Find.Execute("458×5=2290", $true, $false, $false, $false, $false, $true, 1, $false, "844×4=3376", 2) | Out-Null
$d.Content.Find.Execute("225×3=675", $true, $false, $false, $false, $false, $true, 1, $false, "521×8=4168", 2) | Out-Null
$d.Content.Find.Execute("203×4=812", $true, $false, $false, $false, $false, $true, 1, $false, "225×2=450", 2) | Out-Null
$d.Content.Find.Execute("344×9=3096", $true, $false, $false, $false, $false, $true, 1, $false, "938×3=2814", 2) | Out-Null
$d.Content.Find.Execute("146×7=1022", $true, $false, $false, $false, $false, $true, 1, $false, "896×6=5376", 2) | Out-Null
$d.Content.Find.Execute("233×8=1864", $true, $false, $false, $false, $false, $true, 1, $false, "534×4=2136", 2) | Out-Null
$d.Content.Find.Execute("863×5=4315", $true, $false, $false, $false, $false, $true, 1, $false, "384×5=1920", 2) | Out-Null
$d.Content.Find.Execute("334×5=1670", $true, $false, $false, $false, $false, $true, 1, $false, "239×3=717", 2) | Out-Null
$d.Content.Find.Execute("301×9=2709", $true, $false, $false, $false, $false, $true, 1, $false, "605×3=1815", 2) | Out-Null
$d.Content.Find.Execute("151×2=302", $true, $false, $false, $false, $false, $true, 1, $false, "370×5=1850", 2) | Out-Null
$d.Content.Find.Execute("810×3=2430", $true, $false, $false, $false, $false, $true, 1, $false, "376×7=2632", 2) | Out-Null
$d.Content.Find.Execute("652×7=4564", $true, $false, $false, $false, $false, $true, 1, $false, "805×6=4830", 2) | Out-Null
$d.Content.Find.Execute("836×4=3344", $true, $false, $false, $false, $false, $true, 1, $false, "726×7=5082", 2) | Out-Null
$d.Content.Find.Execute("899×5=4495", $true, $false, $false, $false, $false, $true, 1, $false, "683×9=6147", 2) | Out-Null
$d.Content.Find.Execute("110×6=660", $true, $false, $false, $false, $false, $true, 1, $false, "787×4=3148", 2) | Out-Null
$d.Content.Find.Execute("559×6=3354", $true, $false, $false, $false, $false, $true, 1, $false, "136×3=408", 2) | Out-Null
$d.Content.Find.Execute("141×8=1128", $true, $false, $false, $false, $false, $true, 1, $false, "889×9=8001", 2) | Out-Null
$d.Content.Find.Execute("119×8=952", $true, $false, $false, $false, $false, $true, 1, $false, "301×9=2709", 2) | Out-Null
$d.Content.Find.Execute("940×2=1880", $true, $false, $false, $false, $false, $true, 1, $false, "142×8=1136", 2) | Out-Null
$d.Content.Find.Execute("787×9=7083", $true, $false, $false, $false, $false, $true, 1, $false, "230×7=1610", 2) | Out-Null
$d.Content.Find.Execute("297×7=2079", $true, $false, $false, $false, $false, $true, 1, $false, "150×4=600", 2) | Out-Null
$d.Content.Find.Execute("410×3=1230", $true, $false, $false, $false, $false, $true, 1, $false, "113×2=226", 2) | Out-Null
$d.Content.Find.Execute("934×5=4670", $true, $false, $false, $false, $false, $true, 1, $false, "452×3=1356", 2) | Out-Null
$d.Content.Find.Execute("721×3=2163", $true, $false, $false, $false, $false, $true, 1, $false, "628×8=5024", 2) | Out-Null
